$wb = $excel.ActiveWorkbook

# A new handoff occurred for the "b354880c-24dd-4b16-ac5e-af6f3853cae3" file (row 6 of each
# report sheet). Update the "Latest Handoff Date(time)" values to reflect the new handoff,
# while the "cd9fc48a-08cd-4161-afee-beef6704f90c" row (row 7) keeps its previous value.

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("D6").Value = "2016-03-24 22:42:45"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E6").Value = "2016-03-24 22:42:41"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E6").Value = "2016-03-24 22:42:45"
